$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4835600043717818
$ws.Range("C2").Value = 0.2318306736498812
$ws.Range("D2").Value = 0.02719723791074635
$ws.Range("E2").Value = 0.09970402780483312
$ws.Range("F2").Value = 3.785143198636547
$ws.Range("I2").Value = 2.12485044742877
$ws.Range("J2").Value = 0.2065888780660643
$ws.Range("K2").Value = 0.7140377030638376
$ws.Range("M2").Value = 0.2989512149090459

$ws.Range("B3").Value = 0.4641018815264033
$ws.Range("C3").Value = 0.2245224549581053
$ws.Range("D3").Value = 0.02765857281838535
$ws.Range("E3").Value = 0.09945838426490994
$ws.Range("F3").Value = 3.735384344745242
$ws.Range("I3").Value = 2.09663095641686
$ws.Range("J3").Value = 0.2050118552617448
$ws.Range("K3").Value = 0.6874523524119809
$ws.Range("M3").Value = 0.2929052685063702

$ws.Range("B4").Value = 0.4525560959199879
$ws.Range("C4").Value = 0.2201839731732065
$ws.Range("D4").Value = 0.02797041583522741
$ws.Range("E4").Value = 0.09935553210989845
$ws.Range("F4").Value = 3.705944817236286
$ws.Range("I4").Value = 2.079804044332249
$ws.Range("J4").Value = 0.2041161890044663
$ws.Range("K4").Value = 0.671679435236797
$ws.Range("M4").Value = 0.2894000296760098

$ws.Range("B5").Value = 0.4479520672407489
$ws.Range("C5").Value = 0.2184533389796854
$ws.Range("D5").Value = 0.02810470511453644
$ws.Range("E5").Value = 0.09932569896610843
$ws.Range("F5").Value = 3.6942272420699
$ws.Range("I5").Value = 2.073072290808994
$ws.Range("J5").Value = 0.2037694556779783
$ws.Range("K5").Value = 0.6653901955911863
$ws.Range("M5").Value = 0.2880236786987247

$ws.Range("B6").Value = 0.4471936702026653
$ws.Range("C6").Value = 0.218168220890206
$ws.Range("D6").Value = 0.02812743995610489
$ws.Range("E6").Value = 0.0993214752864926
$ws.Range("F6").Value = 3.692298400389745
$ws.Range("I6").Value = 2.071962048272454
$ws.Range("J6").Value = 0.203712983445314
$ws.Range("K6").Value = 0.6643542255510226
$ws.Range("M6").Value = 0.2877982814919413

$ws.Range("B7").Value = 0.4524935955528804
$ws.Range("C7").Value = 0.2201604821864862
$ws.Range("D7").Value = 0.02797219767958836
$ws.Range("E7").Value = 0.09935508083619915
$ws.Range("F7").Value = 3.705785659844494
$ws.Range("I7").Value = 2.079712750517473
$ws.Range("J7").Value = 0.2041114389193481
$ws.Range("K7").Value = 0.6715940562963851
$ws.Range("M7").Value = 0.2893812569107119

$ws.Range("B8").Value = 0.4767674811762106
$ws.Range("C8").Value = 0.2292798642118612
$ws.Range("D8").Value = 0.02735038820386393
$ws.Range("E8").Value = 0.09960938238667794
$ws.Range("F8").Value = 3.767755099136707
$ws.Range("I8").Value = 2.115016363759722
$ws.Range("J8").Value = 0.2060300412326797
$ws.Range("K8").Value = 0.7047567251963471
$ws.Range("M8").Value = 0.2968236136465308

$ws.Range("B9").Value = 0.5275589310404314
$ws.Range("C9").Value = 0.2483486270578226
$ws.Range("D9").Value = 0.02635678120954665
$ws.Range("E9").Value = 0.1004881554124601
$ws.Range("F9").Value = 3.898136612665013
$ws.Range("I9").Value = 2.188234985402772
$ws.Range("J9").Value = 0.2103693824616855
$ws.Range("K9").Value = 0.7741668372294157
$ws.Range("M9").Value = 0.3130613617856284

$ws.Range("B10").Value = 0.5668313854141047
$ws.Range("C10").Value = 0.263090193506514
$ws.Range("D10").Value = 0.02576305880058882
$ws.Range("E10").Value = 0.1013649840910915
$ws.Range("F10").Value = 3.999385476649138
$ws.Range("I10").Value = 2.24449616984802
$ws.Range("J10").Value = 0.2139107220159531
$ws.Range("K10").Value = 0.8278518844747396
$ws.Range("M10").Value = 0.3259962775447818

$ws.Range("B11").Value = 0.5851249543867709
$ws.Range("C11").Value = 0.2699574749775593
$ws.Range("D11").Value = 0.02552227101154614
$ws.Range("E11").Value = 0.1018139898818369
$ws.Range("F11").Value = 4.04664429205701
$ws.Range("I11").Value = 2.270634803216538
$ws.Range("J11").Value = 0.2155988326676948
$ws.Range("K11").Value = 0.852863537412702
$ws.Range("M11").Value = 0.3320998013685852

$ws.Range("B12").Value = 0.592113974435307
$ws.Range("C12").Value = 0.2725812583605318
$ws.Range("D12").Value = 0.02543528149722363
$ws.Range("E12").Value = 0.1019912145990034
$ws.Range("F12").Value = 4.064713357466587
$ws.Range("I12").Value = 2.280611735955389
$ws.Range("J12").Value = 0.2162491875455075
$ws.Range("K12").Value = 0.8624199168862106
$ws.Range("M12").Value = 0.3344426258997686

$ws.Range("B13").Value = 0.5906060215101263
$ws.Range("C13").Value = 0.2720151420447792
$ws.Range("D13").Value = 0.02545383016538239
$ws.Range("E13").Value = 0.1019527262450417
$ws.Range("F13").Value = 4.060814149334021
$ws.Range("I13").Value = 2.278459513805501
$ws.Range("J13").Value = 0.2161086279405922
$ws.Range("K13").Value = 0.8603579961331889
$ws.Range("M13").Value = 0.3339366532589807

$ws.Range("B14").Value = 0.5856987095643831
$ws.Range("C14").Value = 0.2701728677950257
$ws.Range("D14").Value = 0.02551503044909964
$ws.Range("E14").Value = 0.1018284261209743
$ws.Range("F14").Value = 4.048127370800842
$ws.Range("I14").Value = 2.271454029866021
$ws.Range("J14").Value = 0.2156521151448558
$ws.Range("K14").Value = 0.8536480418442522
$ws.Range("M14").Value = 0.3322919145389633

$ws.Range("B15").Value = 0.5827008685245119
$ws.Range("C15").Value = 0.2690474581331443
$ws.Range("D15").Value = 0.02555306262198798
$ws.Range("E15").Value = 0.1017532254551696
$ws.Range("F15").Value = 4.040378923100292
$ws.Range("I15").Value = 2.267173242883047
$ws.Range("J15").Value = 0.2153739343027041
$ws.Range("K15").Value = 0.8495490767083993
$ws.Range("M15").Value = 0.3312885739298181

$ws.Range("B16").Value = 0.5656444808410015
$ws.Range("C16").Value = 0.2626446543784482
$ws.Range("D16").Value = 0.02577938257448409
$ws.Range("E16").Value = 0.1013366481754083
$ws.Range("F16").Value = 3.99632120332825
$ws.Range("I16").Value = 2.242798953209387
$ws.Range("J16").Value = 0.2138019532107904
$ws.Range("K16").Value = 0.8262292010584815
$ws.Range("M16").Value = 0.3256018127380784

$ws.Range("B17").Value = 0.5552906785335097
$ws.Range("C17").Value = 0.2587581188409445
$ws.Range("D17").Value = 0.02592571136280597
$ws.Range("E17").Value = 0.1010939221876725
$ws.Range("F17").Value = 3.969601037238249
$ws.Range("I17").Value = 2.227986029632419
$ws.Range("J17").Value = 0.2128573586942935
$ws.Range("K17").Value = 0.8120744514125704
$ws.Range("M17").Value = 0.3221693587124648

$ws.Range("B18").Value = 0.5493757586973516
$ws.Range("C18").Value = 0.2565378595210746
$ws.Range("D18").Value = 0.02601263489999539
$ws.Range("E18").Value = 0.1009590320358278
$ws.Range("F18").Value = 3.954345237214994
$ws.Range("I18").Value = 2.219517314212865
$ws.Range("J18").Value = 0.2123213127318593
$ws.Range("K18").Value = 0.8039885425693285
$ws.Range("M18").Value = 0.3202157490653477

$ws.Range("B19").Value = 0.547379990945899
$ws.Range("C19").Value = 0.255788720771335
$ws.Range("D19").Value = 0.02604254025548158
$ws.Range("E19").Value = 0.1009141715129473
$ws.Range("F19").Value = 3.949199263870639
$ws.Range("I19").Value = 2.216658748730282
$ws.Range("J19").Value = 0.2121410633103693
$ws.Range("K19").Value = 0.8012603251422661
$ws.Range("M19").Value = 0.3195578363636074

$ws.Range("B20").Value = 0.5563886858509193
$ws.Range("C20").Value = 0.2591702760872749
$ws.Range("D20").Value = 0.02590984901834048
$ws.Range("E20").Value = 0.1011192724370318
$ws.Range("F20").Value = 3.972433752201113
$ws.Range("I20").Value = 2.22955757974691
$ws.Range("J20").Value = 0.2129571609455354
$ws.Range("K20").Value = 0.8135755005304475
$ws.Range("M20").Value = 0.3225326123559213

$ws.Range("B21").Value = 0.5871384326597422
$ws.Range("C21").Value = 0.2707133554069401
$ws.Range("D21").Value = 0.02549694085482912
$ws.Range("E21").Value = 0.101864740885361
$ws.Range("F21").Value = 4.051849080457544
$ws.Range("I21").Value = 2.273509568444368
$ws.Range("J21").Value = 0.2157859025636952
$ws.Range("K21").Value = 0.8556166096258266
$ws.Range("M21").Value = 0.3327741577681635

$ws.Range("B22").Value = 0.6075944963418181
$ws.Range("C22").Value = 0.2783932432776055
$ws.Range("D22").Value = 0.02525150346766836
$ws.Range("E22").Value = 0.1023938854740756
$ws.Range("F22").Value = 4.10476127891792
$ws.Range("I22").Value = 2.302694229392586
$ws.Range("J22").Value = 0.2176993771942293
$ws.Range("K22").Value = 0.8835884842228836
$ws.Range("M22").Value = 0.33965151739676

$ws.Range("B23").Value = 0.5966438140932269
$ws.Range("C23").Value = 0.2742818807032847
$ws.Range("D23").Value = 0.02538027049108393
$ws.Range("E23").Value = 0.1021076379159531
$ws.Range("F23").Value = 4.076428459293908
$ws.Range("I23").Value = 2.287075646971445
$ws.Range("J23").Value = 0.2166721931128279
$ws.Range("K23").Value = 0.8686139677731148
$ws.Range("M23").Value = 0.3359641093923642

$ws.Range("B24").Value = 0.555892159569197
$ws.Range("C24").Value = 0.2589838956929498
$ws.Range("D24").Value = 0.02591701167057892
$ws.Range("E24").Value = 0.101107797084989
$ws.Range("F24").Value = 3.971152751688351
$ws.Range("I24").Value = 2.228846934220329
$ws.Range("J24").Value = 0.2129120185018323
$ws.Range("K24").Value = 0.81289671474687
$ws.Range("M24").Value = 0.3223683238543913

$ws.Range("B25").Value = 0.5134758245855267
$ws.Range("C25").Value = 0.243062171644425
$ws.Range("D25").Value = 0.02660154947604454
$ws.Range("E25").Value = 0.1002097878251327
$ws.Range("F25").Value = 3.861910656092761
$ws.Range("I25").Value = 2.167996707210762
$ws.Range("J25").Value = 0.2091335549173792
$ws.Range("K25").Value = 0.7549187692404189
$ws.Range("M25").Value = 0.3084923671636872
